# Update countries & provincias Spain
#
# The daily COVID country table (sheet "Pais") got refreshed numbers for a
# few countries (Costa de Marfil, Honduras, Islas Virgenes Britanicas,
# Gibraltar). Because the sheet is kept sorted by "Casos totales" (column B)
# descending, those new totals push the affected countries a few rows up,
# shifting the countries that used to occupy those rows down by one. Below
# we just (re)write, row by row, the final country name + stats for every
# row whose content moved or changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Block 1: "Costa de Marfil" gets updated numbers and jumps above
#     "Principado de Andorra" / "Letonia" / "Libano" (rows 88-91). Those
#     three keep their own previous numbers, just shifted one row down.
Set-CountryRow 88 "Costa de Marfil"       688 34 193 489 0  0 6
Set-CountryRow 89 "Principado de Andorra" 682 9  169 480 17 0 33
Set-CountryRow 90 "Letonia"               675 0  57  613 3  0 5
Set-CountryRow 91 "Libano"                663 0  86  556 30 0 21

# --- Block 2: "Honduras" gets updated numbers and jumps above "Nigeria" /
#     "Guinea" (rows 101-103). Those two keep their own previous numbers,
#     just shifted one row down.
Set-CountryRow 101 "Honduras" 442 16 9   392 10 6 41
Set-CountryRow 102 "Nigeria"  442 0  152 277 2  0 13
Set-CountryRow 103 "Guinea"   438 0  49  388 0  0 1

# --- Block 3: "Gibraltar" (row 131) just gets updated numbers, no reorder.
Set-CountryRow 131 "Gibraltar" 132 1 105 27 1 0 0

# --- Block 4: "Islas Virgenes Britanicas" gets updated numbers and jumps
#     above "Bonaire, San Eustaquio y Saba" / "Anguila" (rows 212-214).
#     "San Pedro y Miquelon" and "Yemen" (rows 215-216) also swap order,
#     but they already shared identical stats so no numbers change there.
Set-CountryRow 212 "Islas Virgenes Britanicas"     4 1 2 2 0 0 0
Set-CountryRow 213 "Bonaire, San Eustaquio y Saba" 3 0 0 3 0 0 0
Set-CountryRow 214 "Anguila"                       3 0 1 2 0 0 0
Set-CountryRow 215 "San Pedro y Miquelon"          1 0 0 1 0 0 0
Set-CountryRow 216 "Yemen"                         1 0 0 1 0 0 0
